$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 38, shifting existing rows 38:109 down to 39:110.
$ws.Rows("38:38").Insert()

# Fill in the new row 38 with the new weekly record.
$ws.Cells.Item(38, 1).Value = 5
$ws.Cells.Item(38, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(38, 3).Value = "Maule"
$ws.Cells.Item(38, 4).Value2 = 45203
$ws.Cells.Item(38, 5).Value = 7
$ws.Cells.Item(38, 6).Value = 300000000
$ws.Cells.Item(38, 7).Value = "Espárragos"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 3000
$ws.Cells.Item(38, 11).Value = 1200
$ws.Cells.Item(38, 12).Value = 1200
$ws.Cells.Item(38, 13).Value = 1200
$ws.Cells.Item(38, 14).Value = "`$/kilo"
$ws.Cells.Item(38, 15).Value = "Provincia de Linares"
$ws.Cells.Item(38, 16).Value = 1200
$ws.Cells.Item(38, 17).Value = 1
$ws.Cells.Item(38, 18).Value = "Hortaliza"
